$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new cell N3 - copy format from M3 (empty cell, border style only)
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
# Row 3 height grows from 12.75 (default) to an explicit custom 13.5
$ws.Rows.Item(3).RowHeight = 13.5

# Row 4: new cell N4 - copy format from M4, then set the year value
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2023

# Row 5: new cell N5 - copy format from M5, then set the data value
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 553
